# Adds a 'hole_id' index column to the 'train' worksheet.
# Column A, row 1 gets header "hole_id"; rows 2-38 get the hole id strings
# (previously numeric index values 0..36).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
    "LBU_02_4",
    "LBU_07_03",
    "LBU_01_1",
    "LBU_05_17",
    "LBU_05_01",
    "MHZ_12_03",
    "LBU_05_28",
    "MHZ_12_04",
    "LBU_05_27",
    "LBU_05_21",
    "LBU_05_20",
    "LBU_05_06",
    "MHZ_08_05",
    "LBU_05_14",
    "LBU_05_18",
    "LBU_05_16",
    "LBU_05_29",
    "LBU_01_2",
    "LBU_05_07",
    "LBU_05_09",
    "LBU_05_30",
    "MHZ_12_02",
    "MHZ_08_03",
    "LBU_05_11",
    "LBU_05_08",
    "LBU_05_03",
    "MHZ_12_01",
    "LBU_05_26",
    "LBU_05_23",
    "MHZ_08_02",
    "LBU_07_02",
    "LBU_05_04",
    "LBU_05_05",
    "LBU_05_12",
    "MHZ_08_01",
    "LBU_05_10",
    "LBU_07_01"
)

$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
